$wb = $excel.ActiveWorkbook

# --- Rename sheets (Sheet1/2/3 -> meshing stuff / p4 results / p5 results) ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Name = "meshing stuff"
$ws2.Name = "p4 results"
$ws3.Name = "p5 results"

# --- sheet1: split the big J7:J62 shared formula into the two new rows of results ---
$ws1.Range("J9:J40").Formula = "=H8"
$ws1.Range("J41:J58").Formula = "=H40"

# --- Populate "p4 results" (sheet2) ---
# NB: the order cells are first written controls shared-string allocation
# order, so this follows the same entry order the original author used:
# row-1 headers (skipping the "CMFD iters" column), then the case labels
# down column A, then the row-2 sub-headers, then the "n/a" placeholders,
# and finally the "CMFD iters" header added last.
$ws2.Range("A1").Value = "case"
$ws2.Range("B1").Value = "keff diff"
$ws2.Range("C1").Value = "power"
$ws2.Range("E1").Value = "outers"
$ws2.Range("G1").Value = "runtime"
$ws2.Range("H1").Value = "ratio"

$ws2.Range("A3").Value = "ref"
$ws2.Range("A4").Value = "none"
$ws2.Range("A5").Value = "poly"
$ws2.Range("A6").Value = "subplane"
$ws2.Range("A7").Value = "cpm"

$ws2.Range("C2").Value = "rms"
$ws2.Range("D2").Value = "max"

$ws2.Range("B3").Value = "n/a"
$ws2.Range("C3").Value = "n/a"
$ws2.Range("D3").Value = "n/a"

$ws2.Range("F1").Value = "CMFD iters"

$ws2.Range("E3").Value = 12
$ws2.Range("F3").Value = 517
$ws2.Range("G3").Formula = "=8*60+50"
$ws2.Range("H3").Formula = "=G3/G`$3"

$ws2.Range("B4").Value = -30.280929
$ws2.Range("C4").Value = 0.038360999999999999
$ws2.Range("D4").Value = 0.218083
$ws2.Range("E4").Value = 12
$ws2.Range("F4").Value = 512
$ws2.Range("G4").Formula = "=10*60+4"

$ws2.Range("B5").Value = -7.9185299999999996
$ws2.Range("C5").Value = 0.010293
$ws2.Range("D5").Value = 0.065757999999999997
$ws2.Range("E5").Value = 12
$ws2.Range("F5").Value = 506
$ws2.Range("G5").Formula = "=10*60+8"

$ws2.Range("B6").Value = -7.3639460000000003
$ws2.Range("C6").Value = 0.011278
$ws2.Range("D6").Value = 0.071148000000000003
$ws2.Range("E6").Value = 12
$ws2.Range("F6").Value = 525
$ws2.Range("G6").Formula = "=9*60+58"

$ws2.Range("B7").Value = -1.5902670000000001
$ws2.Range("C7").Value = 0.0054019999999999997
$ws2.Range("D7").Value = 0.049576000000000002
$ws2.Range("E7").Value = 12
$ws2.Range("F7").Value = 526
$ws2.Range("G7").Formula = "=10*60+11"

$ws2.Range("H4:H7").Formula = "=G4/G`$3"

$ws2.Range("C4:C7").NumberFormat = "0.00%"
$ws2.Range("D4:D7").NumberFormat = "0.00%"
$ws2.Range("B4:B7").NumberFormat = "0"

# --- Populate "p5 results" (sheet3) ---
$ws3.Range("A1").Value = "case"
$ws3.Range("B1").Value = "keff diff"
$ws3.Range("C1").Value = "power"
$ws3.Range("E1").Value = "outers"
$ws3.Range("F1").Value = "CMFD iters"
$ws3.Range("G1").Value = "runtime"
$ws3.Range("H1").Value = "ratio"

$ws3.Range("C2").Value = "rms"
$ws3.Range("D2").Value = "max"

$ws3.Range("A3").Value = "ref"
$ws3.Range("B3").Value = "n/a"
$ws3.Range("C3").Value = "n/a"
$ws3.Range("D3").Value = "n/a"
$ws3.Range("E3").Value = 13
$ws3.Range("F3").Value = 445
$ws3.Range("G3").Value = 383.3
$ws3.Range("H3").Formula = "=G3/G`$3"

$ws3.Range("A4").Value = "none"
$ws3.Range("B4").Value = -22
$ws3.Range("C4").Value = 0.0282
$ws3.Range("D4").Value = 0.3055
$ws3.Range("E4").Value = 13
$ws3.Range("F4").Value = 546
$ws3.Range("G4").Value = 439.8

$ws3.Range("A5").Value = "poly"
$ws3.Range("B5").Value = -5
$ws3.Range("C5").Value = 0.0118
$ws3.Range("D5").Value = 0.1178
$ws3.Range("E5").Value = 13
$ws3.Range("F5").Value = 495
$ws3.Range("G5").Value = 407.4

$ws3.Range("A6").Value = "subplane"
$ws3.Range("B6").Value = -5
$ws3.Range("C6").Value = 0.0128
$ws3.Range("D6").Value = 0.1206
$ws3.Range("E6").Value = 13
$ws3.Range("F6").Value = 517
$ws3.Range("G6").Value = 424.3

$ws3.Range("A7").Value = "cpm"
$ws3.Range("B7").Value = -1.0293129999999999
$ws3.Range("C7").Value = 0.0109
$ws3.Range("D7").Value = 0.1217
$ws3.Range("E7").Value = 13
$ws3.Range("F7").Value = 525
$ws3.Range("G7").Formula = "=(60*27+24)/3600*912"

$ws3.Range("H4:H7").Formula = "=G4/G`$3"

$ws3.Range("C4:C7").NumberFormat = "0.00%"
$ws3.Range("D4:D7").NumberFormat = "0.00%"
$ws3.Range("B4:B7").NumberFormat = "0"

# p5 results also picked up an explicit portrait page setup
$ws3.PageSetup.Orientation = 1

# --- Selections / active sheet ---
$ws2.Activate()
$ws2.Range("F8").Select()

$ws3.Range("B4:B7").Select()

$ws1.Range("K8").Select()

$ws2.Activate()

Write-Host "done"
